$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("PLAZAS_TARIFAS")
$ws2 = $wb.Worksheets.Item("PLAZAS_CANALES")

# Fix the duplicated "TORREON " (trailing space) shared-string rows so they
# reuse the canonical "TORREON" string instead. This removes the now-unused
# duplicate entry from the shared-string table and shifts later indices down.
$ws1.Range("A46").Value = "TORREON"
$ws1.Range("A47").Value = "TORREON"
$ws1.Range("A48").Value = "TORREON"
$ws1.Range("A49").Value = "TORREON"

# Data correction: C64 changes from 1 to 24
$ws1.Range("C64").Value = 24

# PLAZAS_CANALES is no longer the active tab; update its stored selection
# first (selecting a range on a sheet activates that sheet as a side effect).
$ws2.Range("B5").Select()

# Update the active sheet / selection to PLAZAS_TARIFAS (this must run last
# so PLAZAS_TARIFAS ends up as the active tab).
$ws1.Activate()
$ws1.Range("C64").Select()

Write-Output "done"
